$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 172: 06-09-2021
$ws.Cells.Item(172, 1).Formula = "=""06-09-2021"""
$ws.Cells.Item(172, 1).Copy()
$ws.Cells.Item(172, 1).PasteSpecial(-4163)
$ws.Cells.Item(172, 2).Value = 1252
$ws.Cells.Item(172, 3).Value = 185
$ws.Cells.Item(172, 4).Value = 316
$ws.Cells.Item(172, 5).Value = 409
$ws.Cells.Item(172, 6).Value = 148
$ws.Cells.Item(172, 7).Value = 194

# New row 173: 07-09-2021
$ws.Cells.Item(173, 1).Formula = "=""07-09-2021"""
$ws.Cells.Item(173, 1).Copy()
$ws.Cells.Item(173, 1).PasteSpecial(-4163)
$ws.Cells.Item(173, 2).Value = 5389
$ws.Cells.Item(173, 3).Value = 442
$ws.Cells.Item(173, 4).Value = 573
$ws.Cells.Item(173, 5).Value = 1532
$ws.Cells.Item(173, 6).Value = 1366
$ws.Cells.Item(173, 7).Value = 1476

$excel.CutCopyMode = $false
